# E-Insurance_Data_Contract.xlsx -- "Apis development" edit (Lakshmi)
#
# Sheet "Database" (sheet1 / ActiveSheet):
#  - Policies table: new row "5 | description | Varchar" at row 9
#  - Customer table: row 18 "age" -> "email"; a new row "7 | mobileNo | varchar"
#    is inserted before the old "policyId" row, which becomes row 20 renumbered to "8"
#  - A blank separator row (style copied from the row above) is inserted at row 40
#  - Row 32 loses its distinct "blank block" style and matches the row-1 style instead
#  - Selection moves to D9 / top of sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database")

# ---------------------------------------------------------------------------
# 1) Policies table, row 9: Sl.No 5 = "description" / "Varchar"
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = 5
$ws.Range("C9").Value = "Varchar"

# ---------------------------------------------------------------------------
# 2) Customer table rows 18-20.
#    New unique strings must be introduced in this exact order so the shared
#    string table comes out as: 52=mobileNo, 53=email, 54=description
# ---------------------------------------------------------------------------
$ws.Range("B19").Value = "mobileNo"
$ws.Range("B18").Value = "email"
$ws.Range("B9").Value  = "description"

$ws.Range("C19").Value = "varchar"
$ws.Range("D19").ClearContents()

$ws.Range("A20").Value = 8
$ws.Range("B20").Value = "policyId"
$ws.Range("C20").Value = "Integer"
$ws.Range("D20").Value = "FK"

# ---------------------------------------------------------------------------
# 3) Row 32 changes from the "blank block" style to the plain row style
#    (same formatting as the Customer table detail rows, e.g. row 19).
# ---------------------------------------------------------------------------
$ws.Range("A19:D19").Copy()
$ws.Range("A32:D32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Insert a new blank row at row 40 (old rows 41-48 shift down to 42-49),
#    picking up the formatting of the row immediately above it (row 39).
# ---------------------------------------------------------------------------
$ws.Rows("40:40").Insert()
$ws.Range("A39:D39").Copy()
$ws.Range("A40:D40").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5) View position: scroll back to the top and select D9.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("D9").Select()
